$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ps = $ws.PageSetup
Write-Host $ps.GetType()
$ps | Get-Member
